# Insert 4 new weekly price rows for "Comercializadora del Agro de Limarí -
# Chirimoya" above the existing row 214, shifting the old rows 214:220 down
# to 218:224 (sheet grows from A1:T220 to A1:T224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 214-220 down by inserting 4 blank rows at 214.
$ws.Range("A214:A217").EntireRow.Insert()

# Common values shared by all four new rows.
$mercadoId = 2
$mercado   = 'Comercializadora del Agro de Limarí'
$region    = 'Coquimbo'
$fecha     = 45267
$codreg    = 4
$tipo      = 'Fruta'
$productoId = 100107
$producto   = 'Otros'
$categoriaId = 100107002
$categoria   = 'Chirimoya'
$variedad    = 'Cultivar IV Región'
$unidad      = '$/bandeja 10 kilos'
$origen      = 'Provincia de Limarí'

# NOTE: this interpreter does not bind PowerShell *named* parameters
# (`-row 214 -calidad 'x'`) correctly, so the helper below is called with
# plain *positional* arguments instead.
function Set-Row {
    param($row, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg, $kgUnidad)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-Row 214 'Especial' 400 16000 17000 16500 1650 10
Set-Row 215 'Primera'  500 12000 13000 12500 1250 10
Set-Row 216 'Segunda'  400 9000  10000 9500  950  10
Set-Row 217 'Tercera'  300 6000  7000  6500  650  10
